# docs: update *1 comment in P.4
# Also refreshes the fixed "Date Placeholder" text (2023/5/28 -> 2023/5/29)
# that is cached on the slide master, every slide layout, and the notes
# master.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Fixed date placeholder text: 2023/5/28 -> 2023/5/29
#    (Slide Master, all 11 Custom Layouts, Notes Master)
# ---------------------------------------------------------------------

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq "2023/5/28") {
                    $shp.TextFrame.TextRange.Text = "2023/5/29"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $lyt = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $lyt.Shapes
}

$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes

# ---------------------------------------------------------------------
# 2. P.4 (Slides.Item(5)) footnote *1 text update
# ---------------------------------------------------------------------

$s = $p.Slides.Item(5)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text.StartsWith("*1 Using ")) {
                $tr = $shp.TextFrame.TextRange

                $r1 = $tr.Characters(1, 9)
                $r1.Text = "*1 Create Mint Account and Metadata Account Using "

                $full = $tr.Text
                $idx = $full.IndexOf(" JavaScript SDK.")
                $start = $idx + 1
                $r3 = $tr.Characters($start, 16)
                $r3.Text = " JavaScript SDK. It's not including upload Metadata JSON and image, verify collection."
            }
        }
    }
}
